$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - interestelar.html
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 0
$ws.Range("G2").Value = 50

# Row 3 - duna.html
$ws.Range("D3").Value = 0
$ws.Range("G3").Value = 50

# Row 4 - matrix.html
$ws.Range("C4").Value = 45
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 45

# Row 5 - mochileiro.html
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 0
$ws.Range("G5").Value = 35
$ws.Range("H5").Value = "Sim"

# Row 6 - blade_runner.html
$ws.Range("C6").Value = 20
$ws.Range("D6").Value = 0
$ws.Range("G6").Value = 10
